# Update the "Pendiente ideal" and "Pendinte real" starting values (day 1)
# from 8 to 9 on the Burndown Chart sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("C6").Value = 9
$ws.Range("C7").Value = 9

# Update sheet view: drop the zoomed-in/scrolled state, select D6 instead of M12.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("D6").Select()
